$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 341.84616
$ws.Range("I28").Value = 203.16667
$ws.Range("J28").Value = 2006
$ws.Range("K28").Value = 203.16667
$ws.Range("L28").Value = 2006
$ws.Range("M28").Value = 281.83333
$ws.Range("N28").Value = -2976

$ws.Range("H51").Value = 4000
$ws.Range("I51").Value = 6000
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 6000
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -5516
$ws.Range("N51").Value = -3968

$ws.Range("H74").Value = 3685.5386
$ws.Range("I74").Value = 804
$ws.Range("J74").Value = 4550
$ws.Range("K74").Value = 804
$ws.Range("L74").Value = 4550
$ws.Range("M74").Value = 132
$ws.Range("N74").Value = -6422

$ws.Range("H76").Value = 3025.1333
$ws.Range("I76").Value = 3032.6086
$ws.Range("J76").Value = 3000.5715
$ws.Range("K76").Value = 3032.6086
$ws.Range("L76").Value = 3000.5715
$ws.Range("M76").Value = -2717.6086
$ws.Range("N76").Value = -3630.5715

$ws.Range("H77").Value = 3685.5386
$ws.Range("I77").Value = 804
$ws.Range("J77").Value = 4550
$ws.Range("K77").Value = 4020
$ws.Range("L77").Value = 22750
$ws.Range("M77").Value = 660
$ws.Range("N77").Value = -32110

$ws.Range("H79").Value = 3025.1333
$ws.Range("I79").Value = 3032.6086
$ws.Range("J79").Value = 3000.5715
$ws.Range("K79").Value = 3032.6086
$ws.Range("L79").Value = 3000.5715
$ws.Range("M79").Value = -1940.6086
$ws.Range("N79").Value = -5184.5715

$ws.Range("H92").Value = 358.9091
$ws.Range("I92").Value = 339.8
$ws.Range("J92").Value = 550
$ws.Range("K92").Value = 339.8
$ws.Range("L92").Value = 550
$ws.Range("M92").Value = 908.2
$ws.Range("N92").Value = -3046

$ws.Range("H107").Value = 839.3913
$ws.Range("I107").Value = 695.7273
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 695.7273
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = 1224.2727
$ws.Range("N107").Value = -7840

$ws.Range("H115").Value = 347
$ws.Range("I115").Value = 347
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1041
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 526
$ws.Range("N115").Value = ""

$ws.Range("H121").Value = 12214.5
$ws.Range("I121").Value = 450
$ws.Range("K121").Value = 1350
$ws.Range("M121").Value = 397

$ws.Range("H132").Value = 2692.325
$ws.Range("I132").Value = 2882.4856
$ws.Range("K132").Value = 8647.4568
$ws.Range("M132").Value = -6117.4568

$ws.Range("H138").Value = 3545.3333
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3545.3333
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 10635.9999
$ws.Range("M138").Value = ""
$ws.Range("N138").Value = -20915.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6337.1807
$ws.Range("I32").Value = 5120.5864
$ws.Range("J32").Value = 11377.357
$ws.Range("K32").Value = 5120.5864
$ws.Range("L32").Value = 11377.357
$ws.Range("M32").Value = -4833.5864
$ws.Range("N32").Value = -11951.357

$ws.Range("H45").Value = 2988.5483
$ws.Range("J45").Value = 2949.0588
$ws.Range("L45").Value = 2949.0588
$ws.Range("N45").Value = -3703.0588

$ws.Range("H63").Value = 2490
$ws.Range("I63").Value = 2548
$ws.Range("J63").Value = 2200
$ws.Range("K63").Value = 2548
$ws.Range("L63").Value = 2200
$ws.Range("M63").Value = -1862
$ws.Range("N63").Value = -3572

$ws.Range("H66").Value = 2490
$ws.Range("I66").Value = 2548
$ws.Range("J66").Value = 2200
$ws.Range("K66").Value = 12740
$ws.Range("L66").Value = 11000
$ws.Range("M66").Value = -9308
$ws.Range("N66").Value = -17864

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = ""

$ws.Range("H97").Value = 38462388
$ws.Range("I97").Value = 600.2105
$ws.Range("J97").Value = 142858670
$ws.Range("K97").Value = 600.2105
$ws.Range("L97").Value = 142858670
$ws.Range("M97").Value = -104.2105
$ws.Range("N97").Value = -142859662

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4719.0625
$ws.Range("I62").Value = 4423
$ws.Range("J62").Value = 6002
$ws.Range("K62").Value = 4423
$ws.Range("L62").Value = 6002
$ws.Range("M62").Value = -3799
$ws.Range("N62").Value = -7250

$ws.Range("H65").Value = 4719.0625
$ws.Range("I65").Value = 4423
$ws.Range("J65").Value = 6002
$ws.Range("K65").Value = 22115
$ws.Range("L65").Value = 30010
$ws.Range("M65").Value = -18995
$ws.Range("N65").Value = -36250

$ws.Range("H107").Value = 1934.625
$ws.Range("I107").Value = 1392.7142
$ws.Range("J107").Value = 2356.111
$ws.Range("K107").Value = 1392.7142
$ws.Range("L107").Value = 2356.111
$ws.Range("M107").Value = 527.2858000000001
$ws.Range("N107").Value = -6196.111

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 65.42856999999999
$ws.Range("I12").Value = 10.5
$ws.Range("K12").Value = 31.5
$ws.Range("M12").Value = 141.5

$ws.Range("H131").Value = 724.0700000000001
$ws.Range("J131").Value = 724.0700000000001
$ws.Range("L131").Value = 2172.21
$ws.Range("N131").Value = -12252.21

$ws.Range("H139").Value = 2093.94
$ws.Range("I139").Value = 1100.9584
$ws.Range("J139").Value = 3010.5386
$ws.Range("K139").Value = 3302.8752
$ws.Range("L139").Value = 9031.6158
$ws.Range("M139").Value = 1837.1248
$ws.Range("N139").Value = -19311.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7731.636
$ws.Range("I70").Value = 3307.6155
$ws.Range("K70").Value = 3307.6155
$ws.Range("M70").Value = -3037.6155

$ws.Range("H73").Value = 7731.636
$ws.Range("I73").Value = 3307.6155
$ws.Range("K73").Value = 3307.6155
$ws.Range("M73").Value = -2371.6155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2811.4211
$ws.Range("I22").Value = 3912.5833
$ws.Range("J22").Value = 923.7143
$ws.Range("K22").Value = 3912.5833
$ws.Range("L22").Value = 923.7143
$ws.Range("M22").Value = -3617.5833
$ws.Range("N22").Value = -1513.7143

$ws.Range("H27").Value = 2811.4211
$ws.Range("I27").Value = 3912.5833
$ws.Range("J27").Value = 923.7143
$ws.Range("K27").Value = 3912.5833
$ws.Range("L27").Value = 923.7143
$ws.Range("M27").Value = -3805.5833
$ws.Range("N27").Value = -1137.7143

$ws.Range("H45").Value = 15666.667
$ws.Range("J45").Value = 15666.667
$ws.Range("L45").Value = 15666.667
$ws.Range("N45").Value = -16480.667

$ws.Range("H46").Value = 2037.0625
$ws.Range("I46").Value = 2135.2144
$ws.Range("K46").Value = 2135.2144
$ws.Range("M46").Value = -1947.2144

$ws.Range("H55").Value = 778.2857
$ws.Range("I55").Value = 1385.7142
$ws.Range("J55").Value = 170.85715
$ws.Range("K55").Value = 1385.7142
$ws.Range("L55").Value = 170.85715
$ws.Range("M55").Value = -1212.7142
$ws.Range("N55").Value = -516.85715

$ws.Range("H93").Value = 1421.8334
$ws.Range("I93").Value = 1506.2
$ws.Range("K93").Value = 1506.2
$ws.Range("M93").Value = -258.2

$ws.Range("H100").Value = 1879.8235
$ws.Range("J100").Value = 2310.4443
$ws.Range("L100").Value = 2310.4443
$ws.Range("N100").Value = -3392.4443

$ws.Range("H104").Value = 19328
$ws.Range("J104").Value = 19328
$ws.Range("L104").Value = 19328
$ws.Range("N104").Value = -26316

$ws.Range("H122").Value = 1637504.4
$ws.Range("I122").Value = 2181417
$ws.Range("K122").Value = 6544251
$ws.Range("M122").Value = -6541801

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4318.273
$ws.Range("I62").Value = 3499.75
$ws.Range("J62").Value = 4786
$ws.Range("K62").Value = 3499.75
$ws.Range("L62").Value = 4786
$ws.Range("M62").Value = -2875.75
$ws.Range("N62").Value = -6034

$ws.Range("H65").Value = 4318.273
$ws.Range("I65").Value = 3499.75
$ws.Range("J65").Value = 4786
$ws.Range("K65").Value = 17498.75
$ws.Range("L65").Value = 23930
$ws.Range("M65").Value = -14378.75
$ws.Range("N65").Value = -30170

$ws.Range("H136").Value = 16924262
$ws.Range("I136").Value = 22441196
$ws.Range("K136").Value = 67323588
$ws.Range("M136").Value = -67321038
